# Applies the "PO Forecast" sheet addition + header renames described by the commit.

$wb = $excel.ActiveWorkbook

$wsWeekly  = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

# --- 1. Rename the "Requested quantity" headers on the existing sheets ---
$wsWeekly.Range("B1").Value  = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2. Add the new "PO Forecast" worksheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Reuse the existing header / date cell styles (bold header, date-formatted column A)
# by copying formats only, so the new sheet shares the same style indices.
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A39").PasteSpecial(-4122)

# --- 3. Header row ---
$wsForecast.Cells.Item(1,1).Value = "ds"
$wsForecast.Cells.Item(1,2).Value = "PO_Forecast"
$wsForecast.Cells.Item(1,3).Value = "yhat_lower"
$wsForecast.Cells.Item(1,4).Value = "yhat_upper"

# --- 4. Forecast data rows (ds, PO_Forecast, yhat_lower, yhat_upper) ---
$rows = @(
    @(2, 45333.99999999999, 24, -8.702351908703118, 53.54625103072828),
    @(3, 45340.99999999999, 23, -7.703764037613189, 55.45000498496821),
    @(4, 45354.99999999999, 23, -9.312305944756703, 56.36319995780222),
    @(5, 45361.99999999999, 23, -8.687269704633515, 56.99103575672824),
    @(6, 45368.99999999999, 23, -6.979174899783488, 52.78371162970051),
    @(7, 45375.99999999999, 22, -9.640833120776811, 55.42867717881278),
    @(8, 45382.99999999999, 22, -10.40425551204973, 53.01507944241717),
    @(9, 45389.99999999999, 22, -9.730372022917676, 52.15849609490646),
    @(10, 45396.99999999999, 22, -11.76990673652718, 54.53651567408031),
    @(11, 45403.99999999999, 22, -9.957232805792081, 51.42719444583287),
    @(12, 45410.99999999999, 21, -9.433269939160844, 51.78292819572462),
    @(13, 45417.99999999999, 21, -9.988403178793812, 54.54058804869855),
    @(14, 45424.99999999999, 21, -9.064358876353568, 53.45438779218522),
    @(15, 45431.99999999999, 21, -10.41134416905788, 52.56871093994184),
    @(16, 45438.99999999999, 21, -11.69575392034239, 50.59244309179051),
    @(17, 45445.99999999999, 20, -9.359812581542608, 53.0664799099566),
    @(18, 45452.99999999999, 20, -10.74124634181103, 53.66298033337151),
    @(19, 45473.99999999999, 20, -13.05638603659709, 52.10098540021083),
    @(20, 45480.99999999999, 20, -11.94949477514512, 52.51871991587377),
    @(21, 45487.99999999999, 19, -11.31702757983164, 49.69940284469172),
    @(22, 45529.99999999999, 18, -14.35144734936789, 48.97361033084076),
    @(23, 45536.99999999999, 18, -15.61534074494302, 50.7688435967352),
    @(24, 45550.99999999999, 18, -14.62261013370857, 49.62302589034244),
    @(25, 45564.99999999999, 17, -15.05456846928768, 48.48861445550723),
    @(26, 45571.99999999999, 17, -14.78341716282825, 45.28271198308686),
    @(27, 45578.99999999999, 17, -14.59024503388187, 49.11293695077377),
    @(28, 45585.99999999999, 17, -16.28283779367754, 49.78195458869033),
    @(29, 45613.99999999999, 16, -14.41352262753584, 48.94382763788506),
    @(30, 45634.99999999999, 15, -18.11418694952163, 48.53809014554446),
    @(31, 45641.99999999999, 15, -17.61796307045955, 48.39190203070866),
    @(32, 45648.99999999999, 15, -17.28551684998975, 48.13925182800165),
    @(33, 45655.99999999999, 15, -17.87331252803848, 45.22414670364693),
    @(34, 45662.99999999999, 14, -16.77102439295157, 46.38610853051541),
    @(35, 45669.99999999999, 14, -16.13980796440663, 45.06513883372509),
    @(36, 45676.99999999999, 14, -17.95689344129537, 45.87620067042811),
    @(37, 45683.99999999999, 14, -18.69516030085298, 45.28135640654769),
    @(38, 45690.99999999999, 14, -16.56251250410292, 44.90298087705142),
    @(39, 45697.99999999999, 14, -16.29530395307235, 45.70255466294547)
)

foreach ($row in $rows) {
    $r = $row[0]
    $wsForecast.Cells.Item($r, 1).Value = $row[1]
    $wsForecast.Cells.Item($r, 2).Value = $row[2]
    $wsForecast.Cells.Item($r, 3).Value = $row[3]
    $wsForecast.Cells.Item($r, 4).Value = $row[4]
}
